$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.460.78'
$ws.Range('E2').Value = '  -0.71%  '
$ws.Range('D3').Value = '2.376.83'
$ws.Range('E3').Value = '  +5.91%  '
$ws.Range('E4').Value = '  -0.57%  '
$ws.Range('D5').Value = "'0.653"
$ws.Range('E5').Value = '  +2.04%  '
$ws.Range('D6').Value = "'232.62"
$ws.Range('E6').Value = '  +1.19%  '
$ws.Range('D7').Value = "'68.17"
$ws.Range('E7').Value = '  +6.30%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = "'0.459"
$ws.Range('E9').Value = '  +2.89%  '
$ws.Range('D10').Value = "'0.0949"
$ws.Range('E10').Value = '  -2.41%  '
$ws.Range('D11').Value = "'56.94"
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('D12').Value = "'26.48"
$ws.Range('E12').Value = '  -0.42%  '
$ws.Range('D13').Value = '2.727.68'
$ws.Range('E13').Value = '  +5.53%  '
$ws.Range('E14').Value = '  -1.28%  '
$ws.Range('D15').Value = "'15.58"
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('D16').Value = "'6.29"
$ws.Range('E16').Value = '  +3.79%  '
$ws.Range('D17').Value = "'0.844"
$ws.Range('E17').Value = '  +1.84%  '
$ws.Range('D18').Value = '2.379.01'
$ws.Range('E18').Value = '  +5.15%  '
$ws.Range('D19').Value = '43.457.21'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').Value = '0.0₃0982'
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('D21').Value = "'6.29"
$ws.Range('E21').Value = '  +4.51%  '
$ws.Range('D22').Value = "'73.72"
$ws.Range('E22').Value = '  +0.91%  '
$ws.Range('D23').Value = "'248.77"
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('D24').Value = "'3.90"
$ws.Range('E24').Value = '  +17.85%  '
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('E26').Value = '  +1.20%  '
$ws.Range('E27').Value = '  -1.26%  '
$ws.Range('E28').Value = '  +0.32%  '
$ws.Range('D29').Value = "'22.45"
$ws.Range('E29').Value = '  +8.16%  '
$ws.Range('D30').Value = "'175.19"
$ws.Range('E30').Value = '  +2.65%  '
$ws.Range('D31').Value = "'1.51"
$ws.Range('E31').Value = '  +9.67%  '
$ws.Range('D32').Value = "'0.130"
$ws.Range('E32').Value = '  -5.97%  '
$ws.Range('E33').Value = '  +1.28%  '
$ws.Range('D34').Value = "'5.00"
$ws.Range('E34').Value = '  +5.47%  '
$ws.Range('D35').Value = "'0.0696"
$ws.Range('E35').Value = '  -0.50%  '
$ws.Range('D36').Value = "'5.04"
$ws.Range('E36').Value = '  +3.69%  '
$ws.Range('D37').Value = "'2.53"
$ws.Range('E37').Value = '  +11.30%  '
$ws.Range('D38').Value = "'6.51"
$ws.Range('E38').Value = '  +1.36%  '
$ws.Range('D39').Value = "'3.64"
$ws.Range('E39').Value = '  -1.76%  '
$ws.Range('E40').Value = '  -1.96%  '
$ws.Range('E41').Value = '  +10.58%  '
$ws.Range('D42').Value = "'1.00"
$ws.Range('E42').Value = '  -0.14%  '
$ws.Range('D43').Value = "'17.88"
$ws.Range('E43').Value = '  +3.87%  '
$ws.Range('E44').Value = '  +9.49%  '
$ws.Range('D45').Value = "'99.60"
$ws.Range('E45').Value = '  +2.72%  '
$ws.Range('E46').Value = '  +1.77%  '
$ws.Range('D47').Value = "'0.0952"
$ws.Range('E47').Value = '  -1.13%  '
$ws.Range('D48').Value = "'4.37"
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('D49').Value = '1.448.25'
$ws.Range('E49').Value = '  +1.26%  '
$ws.Range('D50').Value = '2.601.24'
$ws.Range('E50').Value = '  +5.77%  '
$ws.Range('E51').Value = '  -6.92%  '

$clearCells = @('D5','D6','D7','D9','D10','D11','D12','D15','D16','D17','D21','D22','D23','D24','D29','D30','D31','D32','D34','D35','D36','D37','D38','D39','D42','D43','D45','D47','D48')
foreach ($c in $clearCells) { $ws.Range($c).ClearFormats() }
